$d = $word.ActiveDocument

function FindParaIndex($text, $matchCase) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $matchCase, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $text"
    }
    return $rng.Paragraphs.Item(1).Index
}

# ---------------------------------------------------------------------------
# Step 1: apply "No Spacing" style to the duplicated Field (1)/(2) block that
# currently still carries the default Normal style (2nd "Test Case: Field (1)"
# through the 2nd "Input: " paragraph).
# ---------------------------------------------------------------------------
$idxTestObjective2 = FindParaIndex "Input contains one compulsory field, which is not specific" $true
$idxTestCase2 = $idxTestObjective2 - 2
$d.Paragraphs.Item($idxTestCase2).Style = "No Spacing"
$d.Paragraphs.Item($idxTestCase2 + 1).Style = "No Spacing"
$d.Paragraphs.Item($idxTestCase2 + 2).Style = "No Spacing"
$d.Paragraphs.Item($idxTestCase2 + 3).Style = "No Spacing"

Write-Output "Step1 done idx=$idxTestCase2"

# ---------------------------------------------------------------------------
# Step 2: collapse the run-fragmented "*Data is stored inside Field (2).txt…"
# paragraph (2nd "Test Files" block) down to a single plain run, removing the
# italic formatting that had been applied to "Test Files".
# ---------------------------------------------------------------------------
$idxField2Data = FindParaIndex "Field (2" $true
$pField2 = $d.Paragraphs.Item($idxField2Data)
$rField2 = $pField2.Range
$rField2.MoveEnd(1, -1) | Out-Null
$rField2.Delete()
$rField2.InsertAfter("*Data is stored inside " + [char]8220 + "Field (2).txt" + [char]8221 + " inside Test Files folder.")

Write-Output "Step2 done idx=$idxField2Data"

# ---------------------------------------------------------------------------
# Step 3: the lone empty paragraph right after the Field (2) "Expected
# behavior" paragraph (and right before the "[8] - Title" block) is expanded
# into four new paragraphs that document the new "Field (3)" test case intro.
# ---------------------------------------------------------------------------
$idxExpected2 = FindParaIndex "A new variable is declared." $true
$idxEmptyBeforeTitle = $idxExpected2 + 1
$pEmpty = $d.Paragraphs.Item($idxEmptyBeforeTitle)
$rEmpty = $pEmpty.Range

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newField3Intro = (
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t>Test Case: Field (3)</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr><w:r><w:t xml:space="preserve">Date designed: </w:t></w:r><w:r><w:t>11 June 2012</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Test Objective: A valid entry of typical field inside a segment. </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>Input contain</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> one non-compulsory field, which is specific from a list. List is given in the </w:t></w:r>' +
    '<w:r><w:t>Field segment.</w:t></w:r>' +
  '</w:p>'
)

$rEmpty.InsertXML($xmlHeader + $newField3Intro + $xmlFooter) | Out-Null

Write-Output "Step3 done idx=$idxEmptyBeforeTitle"

# ---------------------------------------------------------------------------
# Step 4: append the new "Field (3)" wrap-up paragraphs right after the final
# "-----" separator (the one that closes out the "[8] - Title" block), and
# relocate the _GoBack bookmark so its start sits in the "[8] - Title"
# paragraph and its end sits just before the newly appended paragraphs.
# ---------------------------------------------------------------------------
$idxProf = FindParaIndex "100568 - Prof" $true
$idxFinalDashes = $idxProf + 1
$pFinalDashes = $d.Paragraphs.Item($idxFinalDashes)
Write-Output "finalDashes text=[$($pFinalDashes.Range.Text)]"

$idxTitle = FindParaIndex "[8] - Title" $true
$pTitle = $d.Paragraphs.Item($idxTitle)
Write-Output "title text=[$($pTitle.Range.Text)]"

# Collapsed range right after the final "-----" paragraph's mark — this is
# where the new paragraphs get spliced in.
$insertPoint = $d.Range($pFinalDashes.Range.End, $pFinalDashes.Range.End)

$newField3Outro = (
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:r><w:t>*Data is stored inside ' + [char]8220 + 'Field (3</w:t></w:r>' +
    '<w:r><w:t>).txt' + [char]8221 + ' inside Test Files folder.</w:t></w:r>' +
  '</w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Expected behavior: Code will run without error. </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">A new mapping array will be formed. </w:t></w:r>' +
    '<w:r><w:t>A new variable is declared.</w:t></w:r>' +
  '</w:p>'
)

$insertPoint.InsertXML($xmlHeader + $newField3Outro + $xmlFooter) | Out-Null

Write-Output "Step4 done idxFinalDashes=$idxFinalDashes idxTitle=$idxTitle"
